$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update financial figures for rows 2-6 (existing cells get new values) ---
# Row 2
$ws.Range("D2").Value = 4318
$ws.Range("E2").Value = 151
$ws.Range("F2").Value = 240
$ws.Range("G2").Value = 319
$ws.Range("H2").Value = 247
$ws.Range("I2").Value = 236
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 8183
$ws.Range("L2").Value = 1522
$ws.Range("M2").Value = 6661
$ws.Range("N2").Value = 6507
$ws.Range("O2").Value = 154
$ws.Range("P2").Value = 120
$ws.Range("Q2").Value = 238
$ws.Range("R2").Value = -445
$ws.Range("S2").Value = 213
$ws.Range("T2").Value = 242
$ws.Range("U2").Value = -4
$ws.Range("V2").Value = 397
$ws.Range("W2").Value = 3.49
$ws.Range("X2").Value = 5.73
$ws.Range("Y2").Value = 3.67
$ws.Range("Z2").Value = 3.1
$ws.Range("AA2").Value = 22.85
$ws.Range("AB2").Value = 5321.22
$ws.Range("AC2").Value = 9817
$ws.Range("AD2").Value = 17.52
$ws.Range("AE2").Value = 277593
$ws.Range("AF2").Value = 0.62
$ws.Range("AG2").Value = 1750
$ws.Range("AH2").Value = 1.02
$ws.Range("AI2").Value = 17.41
$ws.Range("AJ2").Value = 2400000
# Row 3
$ws.Range("D3").Value = 4547
$ws.Range("E3").Value = 177
$ws.Range("F3").Value = 255
$ws.Range("G3").Value = 281
$ws.Range("H3").Value = 202
$ws.Range("I3").Value = 193
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 8895
$ws.Range("L3").Value = 2106
$ws.Range("M3").Value = 6789
$ws.Range("N3").Value = 6627
$ws.Range("O3").Value = 162
$ws.Range("P3").Value = 120
$ws.Range("Q3").Value = 180
$ws.Range("R3").Value = -829
$ws.Range("S3").Value = 581
$ws.Range("T3").Value = 607
$ws.Range("U3").Value = -426
$ws.Range("V3").Value = 1059
$ws.Range("W3").Value = 3.9
$ws.Range("X3").Value = 4.45
$ws.Range("Y3").Value = 2.93
$ws.Range("Z3").Value = 2.37
$ws.Range("AA3").Value = 31.03
$ws.Range("AB3").Value = 5446.31
$ws.Range("AC3").Value = 8026
$ws.Range("AD3").Value = 19.19
$ws.Range("AE3").Value = 283294
$ws.Range("AF3").Value = 0.54
$ws.Range("AG3").Value = 2500
$ws.Range("AH3").Value = 1.62
$ws.Range("AI3").Value = 30.36
$ws.Range("AJ3").Value = 2400000
# Row 4
$ws.Range("D4").Value = 4783
$ws.Range("E4").Value = 206
$ws.Range("F4").Value = 277
$ws.Range("G4").Value = 281
$ws.Range("H4").Value = 207
$ws.Range("I4").Value = 198
$ws.Range("J4").Value = 9
$ws.Range("K4").Value = 9048
$ws.Range("L4").Value = 2106
$ws.Range("M4").Value = 6942
$ws.Range("N4").Value = 6773
$ws.Range("O4").Value = 169
$ws.Range("P4").Value = 120
$ws.Range("Q4").Value = 392
$ws.Range("R4").Value = -128
$ws.Range("S4").Value = -123
$ws.Range("T4").Value = 159
$ws.Range("U4").Value = 232
$ws.Range("V4").Value = 1030
$ws.Range("W4").Value = 4.3
$ws.Range("X4").Value = 4.33
$ws.Range("Y4").Value = 2.96
$ws.Range("Z4").Value = 2.31
$ws.Range("AA4").Value = 30.34
$ws.Range("AB4").Value = 5561.93
$ws.Range("AC4").Value = 8265
$ws.Range("AD4").Value = 14.94
$ws.Range("AE4").Value = 289667
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 2800
$ws.Range("AH4").Value = 2.27
$ws.Range("AI4").Value = 33
$ws.Range("AJ4").Value = 2400000
# Row 5
$ws.Range("D5").Value = 4925
$ws.Range("E5").Value = 145
$ws.Range("F5").Value = 240
$ws.Range("G5").Value = 289
$ws.Range("H5").Value = 170
$ws.Range("I5").Value = 160
$ws.Range("J5").Value = 11
$ws.Range("K5").Value = 8959
$ws.Range("L5").Value = 1910
$ws.Range("M5").Value = 7049
$ws.Range("N5").Value = 6870
$ws.Range("O5").Value = 179
$ws.Range("P5").Value = 120
$ws.Range("Q5").Value = 370
$ws.Range("R5").Value = -207
$ws.Range("S5").Value = -179
$ws.Range("T5").Value = 135
$ws.Range("U5").Value = 235
$ws.Range("V5").Value = 842
$ws.Range("W5").Value = 2.94
$ws.Range("X5").Value = 3.46
$ws.Range("Y5").Value = 2.34
$ws.Range("Z5").Value = 1.89
$ws.Range("AA5").Value = 27.1
$ws.Range("AB5").Value = 5640.09
$ws.Range("AC5").Value = 6655
$ws.Range("AD5").Value = 16.53
$ws.Range("AE5").Value = 296290
$ws.Range("AF5").Value = 0.37
$ws.Range("AG5").Value = 1500
$ws.Range("AH5").Value = 1.36
$ws.Range("AI5").Value = 21.78
$ws.Range("AJ5").Value = 2400000
# Row 6
$ws.Range("D6").Value = 4923
$ws.Range("E6").Value = 190
$ws.Range("F6").Value = 299
$ws.Range("G6").Value = 306
$ws.Range("H6").Value = 223
$ws.Range("I6").Value = 213
$ws.Range("K6").Value = 9373
$ws.Range("L6").Value = 2141
$ws.Range("M6").Value = 7232
$ws.Range("N6").Value = 7045
$ws.Range("P6").Value = 120
$ws.Range("Q6").Value = 75
$ws.Range("R6").Value = -315
$ws.Range("S6").Value = 165
$ws.Range("T6").Value = 266
$ws.Range("U6").Value = -191
$ws.Range("V6").Value = 1074
$ws.Range("W6").Value = 3.86
$ws.Range("X6").Value = 4.54
$ws.Range("Y6").Value = 3.06
$ws.Range("Z6").Value = 2.44
$ws.Range("AA6").Value = 29.6
$ws.Range("AB6").Value = 5818.04
$ws.Range("AC6").Value = 8865
$ws.Range("AD6").Value = 10.83
$ws.Range("AE6").Value = 306851
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 2.08
$ws.Range("AI6").Value = 21.58
$ws.Range("AJ6").Value = 2400000

# --- Rows 7-9: clear all data cells except A (row #), B (연간), C (year label) ---
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
